$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2023-04-30 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-01 Monday", 2) | Out-Null

# Update each multiplication-problem cell in the table by (row, column) index
# to avoid ambiguity from duplicate / overlapping text values across cells.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50×45="
$t.Cell(1, 2).Range.Text = "58×76="
$t.Cell(1, 3).Range.Text = "18×100="
$t.Cell(1, 4).Range.Text = "16×74="
$t.Cell(1, 5).Range.Text = "43×52="
$t.Cell(2, 1).Range.Text = "77×70="
$t.Cell(2, 2).Range.Text = "93×65="
$t.Cell(2, 3).Range.Text = "80×10="
$t.Cell(2, 4).Range.Text = "95×61="
$t.Cell(2, 5).Range.Text = "83×76="
$t.Cell(3, 1).Range.Text = "40×41="
$t.Cell(3, 2).Range.Text = "51×56="
$t.Cell(3, 3).Range.Text = "10×78="
$t.Cell(3, 4).Range.Text = "99×18="
$t.Cell(3, 5).Range.Text = "12×82="
$t.Cell(4, 1).Range.Text = "92×10="
$t.Cell(4, 2).Range.Text = "55×44="
$t.Cell(4, 3).Range.Text = "97×32="
$t.Cell(4, 4).Range.Text = "90×98="
$t.Cell(4, 5).Range.Text = "47×25="
$t.Cell(5, 1).Range.Text = "36×71="
$t.Cell(5, 2).Range.Text = "36×23="
$t.Cell(5, 3).Range.Text = "75×40="
$t.Cell(5, 4).Range.Text = "11×18="
$t.Cell(5, 5).Range.Text = "26×59="
$t.Cell(6, 1).Range.Text = "74×31="
$t.Cell(6, 2).Range.Text = "63×15="
$t.Cell(6, 3).Range.Text = "94×98="
$t.Cell(6, 4).Range.Text = "85×15="
$t.Cell(6, 5).Range.Text = "92×77="
$t.Cell(7, 1).Range.Text = "16×47="
$t.Cell(7, 2).Range.Text = "23×13="
$t.Cell(7, 3).Range.Text = "49×22="
$t.Cell(7, 4).Range.Text = "55×40="
$t.Cell(7, 5).Range.Text = "81×60="
$t.Cell(8, 1).Range.Text = "30×81="
$t.Cell(8, 2).Range.Text = "33×15="
$t.Cell(8, 3).Range.Text = "72×73="
$t.Cell(8, 4).Range.Text = "94×60="
$t.Cell(8, 5).Range.Text = "11×39="
$t.Cell(9, 1).Range.Text = "52×61="
$t.Cell(9, 2).Range.Text = "89×78="
$t.Cell(9, 3).Range.Text = "62×80="
$t.Cell(9, 4).Range.Text = "49×64="
$t.Cell(9, 5).Range.Text = "26×15="
$t.Cell(10, 1).Range.Text = "60×37="
$t.Cell(10, 2).Range.Text = "38×72="
$t.Cell(10, 3).Range.Text = "85×28="
$t.Cell(10, 4).Range.Text = "90×43="
$t.Cell(10, 5).Range.Text = "94×44="
$t.Cell(11, 1).Range.Text = "41×26="
$t.Cell(11, 2).Range.Text = "65×43="
$t.Cell(11, 3).Range.Text = "54×56="
$t.Cell(11, 4).Range.Text = "58×89="
$t.Cell(11, 5).Range.Text = "17×13="
$t.Cell(12, 1).Range.Text = "19×61="
$t.Cell(12, 2).Range.Text = "28×58="
$t.Cell(12, 3).Range.Text = "12×60="
$t.Cell(12, 4).Range.Text = "83×15="
$t.Cell(12, 5).Range.Text = "15×78="
$t.Cell(13, 1).Range.Text = "74×21="
$t.Cell(13, 2).Range.Text = "56×95="
$t.Cell(13, 3).Range.Text = "93×37="
$t.Cell(13, 4).Range.Text = "94×89="
$t.Cell(13, 5).Range.Text = "56×33="
$t.Cell(14, 1).Range.Text = "63×72="
$t.Cell(14, 2).Range.Text = "67×20="
$t.Cell(14, 3).Range.Text = "76×51="
$t.Cell(14, 4).Range.Text = "43×20="
$t.Cell(14, 5).Range.Text = "58×42="
$t.Cell(15, 1).Range.Text = "38×33="
$t.Cell(15, 2).Range.Text = "35×85="
$t.Cell(15, 3).Range.Text = "58×10="
$t.Cell(15, 4).Range.Text = "67×29="
$t.Cell(15, 5).Range.Text = "32×44="
$t.Cell(16, 1).Range.Text = "44×29="
$t.Cell(16, 2).Range.Text = "67×63="
$t.Cell(16, 3).Range.Text = "72×85="
$t.Cell(16, 4).Range.Text = "85×22="
$t.Cell(16, 5).Range.Text = "93×68="
$t.Cell(17, 1).Range.Text = "77×19="
$t.Cell(17, 2).Range.Text = "83×32="
$t.Cell(17, 3).Range.Text = "29×84="
$t.Cell(17, 4).Range.Text = "88×56="
$t.Cell(17, 5).Range.Text = "94×53="
$t.Cell(18, 1).Range.Text = "92×17="
$t.Cell(18, 2).Range.Text = "76×70="
$t.Cell(18, 3).Range.Text = "46×90="
$t.Cell(18, 4).Range.Text = "32×90="
$t.Cell(18, 5).Range.Text = "100×56="
$t.Cell(19, 1).Range.Text = "96×53="
$t.Cell(19, 2).Range.Text = "38×53="
$t.Cell(19, 3).Range.Text = "76×83="
$t.Cell(19, 4).Range.Text = "31×76="
$t.Cell(19, 5).Range.Text = "54×74="
$t.Cell(20, 1).Range.Text = "17×94="
$t.Cell(20, 2).Range.Text = "77×14="
$t.Cell(20, 3).Range.Text = "66×49="
$t.Cell(20, 4).Range.Text = "57×17="
$t.Cell(20, 5).Range.Text = "52×39="
